$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item(1).Name = "GNG_TO-1650291278980096"
$wb.Worksheets.Item(2).Name = "NB_TO-1650291283268005"
$wb.Worksheets.Item(3).Name = "RS_TO-16502912832700047"
$wb.Worksheets.Item(4).Name = "TOL_TO-16502912833203893"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16502912833834398"

# Sheet 1 (GNG)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16502912789325266.csv"
$ws1.Range("B3").Value = "GNG_stims-16502912789475658.csv"
$ws1.Range("B4").Value = "go_stims-16502912789499106.csv"
$ws1.Range("B5").Value = "GNG_stims-1650291278979094.csv"

# Sheet 2 (NB)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16502912805400288.csv"
$ws2.Range("B3").Value = "ZB-match_5-16502912796173196.csv"
$ws2.Range("B4").Value = "TB-1650291283252908.csv"
$ws2.Range("B5").Value = "ZB-match_4-16502912792958722.csv"
$ws2.Range("B6").Value = "OB-16502912803278086.csv"
$ws2.Range("B7").Value = "TB-1650291281917672.csv"
$ws2.Range("B8").Value = "TB-16502912824346256.csv"
$ws2.Range("B9").Value = "ZB-match_5-1650291279779899.csv"
$ws2.Range("B10").Value = "OB-16502912812629676.csv"

# Sheet 3 (RS)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# Sheet 4 (TOL)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16502912832834225.csv"
$ws4.Range("B3").Value = "ZM_stims-16502912832720032.csv"
$ws4.Range("B4").Value = "MM_stims-16502912833034446.csv"
$ws4.Range("B5").Value = "ZM_stims-16502912832844234.csv"
$ws4.Range("B6").Value = "MM_stims-16502912833193905.csv"
$ws4.Range("B7").Value = "ZM_stims-16502912833044484.csv"

# Sheet 5 (vSAT)
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-1650291283338133.csv"
$ws5.Range("B3").Value = "vSAT_stims-16502912833693855.csv"
$ws5.Range("B4").Value = "SAT_stims-16502912833258204.csv"
$ws5.Range("B5").Value = "vSAT_stims-16502912833539646.csv"
